# "add the Scrollable to the dataTable"
# Extend Sheet1's data table with three more repeated blocks of rows (the
# original had a header row + 7 data rows; we append 17 more data rows),
# widen column B to fit the new, longer department names, and make Sheet1
# (rather than Sheet2) the active sheet/selection again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data rows (9-25) -------------------------------------------------
# Rows 9-13 are entered in a specific cell order so that the handful of
# brand-new strings land in the shared-strings table in ascending order
# (matches how the source workbook's shared-strings table was built up).
$ws.Range("A9").Value = "季昌明"
$ws.Range("B9").Value = "汉东大学"
$ws.Range("C9").Value = "检查"
$ws.Range("D9").Value = "检察长"

$ws.Range("A10").Value = "陆亦可"
$ws.Range("B10").Value = "汉东大学"
$ws.Range("C10").Value = "侦查"
$ws.Range("D10").Value = "侦查处长"

$ws.Range("A11").Value = "林华华"
$ws.Range("C11").Value = "检查"
$ws.Range("D11").Value = "美女科长"
$ws.Range("A12").Value = "周正"
$ws.Range("B11").Value = "汉东大学政法系2009级"
$ws.Range("D12").Value = "科长"
$ws.Range("B12").Value = "汉东大学政法系2009级"
$ws.Range("C12").Value = "检查"

$ws.Range("A13").Value = "高玉亮"
$ws.Range("B13").Value = "汉东大学政法系1990级"
$ws.Range("C13").Value = "政法"
$ws.Range("D13").Value = "政法委书记"

$ws.Range("A14").Value = "季昌明"
$ws.Range("B14").Value = "汉东大学"
$ws.Range("C14").Value = "检查"
$ws.Range("D14").Value = "检察长"

$ws.Range("A15").Value = "陆亦可"
$ws.Range("B15").Value = "汉东大学"
$ws.Range("C15").Value = "侦查"
$ws.Range("D15").Value = "侦查处长"

$ws.Range("A16").Value = "林华华"
$ws.Range("B16").Value = "汉东大学政法系2009级"
$ws.Range("C16").Value = "检查"
$ws.Range("D16").Value = "美女科长"

$ws.Range("A17").Value = "周正"
$ws.Range("B17").Value = "汉东大学政法系2009级"
$ws.Range("C17").Value = "检查"
$ws.Range("D17").Value = "科长"

$ws.Range("A18").Value = "高玉亮"
$ws.Range("B18").Value = "汉东大学政法系1990级"
$ws.Range("C18").Value = "政法"
$ws.Range("D18").Value = "政法委书记"

$ws.Range("A19").Value = "侯亮平"
$ws.Range("B19").Value = "中国传媒大学"
$ws.Range("C19").Value = "跳舞"
$ws.Range("D19").Value = "舞蹈家"

$ws.Range("A20").Value = "陈海"
$ws.Range("B20").Value = "中国传媒大学"
$ws.Range("C20").Value = "武术"
$ws.Range("D20").Value = "武术家"

$ws.Range("A21").Value = "沙瑞金"
$ws.Range("B21").Value = "中国传媒大学"
$ws.Range("C21").Value = "表演"
$ws.Range("D21").Value = "一级演员"

$ws.Range("A22").Value = "季昌明"
$ws.Range("B22").Value = "汉东大学"
$ws.Range("C22").Value = "检查"
$ws.Range("D22").Value = "检察长"

$ws.Range("A23").Value = "陆亦可"
$ws.Range("B23").Value = "汉东大学"
$ws.Range("C23").Value = "侦查"
$ws.Range("D23").Value = "侦查处长"

$ws.Range("A24").Value = "林华华"
$ws.Range("B24").Value = "汉东大学政法系2009级"
$ws.Range("C24").Value = "检查"
$ws.Range("D24").Value = "美女科长"

$ws.Range("A25").Value = "周正"
$ws.Range("B25").Value = "汉东大学政法系2009级"
$ws.Range("C25").Value = "检查"
$ws.Range("D25").Value = "科长"

# --- Column B is now much wider (longer department names) -----------------
$ws.Columns.Item(2).ColumnWidth = 20.71

# --- Sheet1 becomes the active sheet / selection again (Sheet2's own
#     E7 selection is untouched, so it is left exactly as it was) ----------
$ws.Activate() | Out-Null
$ws.Range("E17").Select() | Out-Null
